$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Pipeline(steps=[('scaler', RobustScaler()), ('selector', None),
                ('model',
                 LGBMClassifier(boosting_type='dart', class_weight='balanced',
                                colsample_bytree=0.7, learning_rate=0.01,
                                max_depth=7, min_child_samples=1, num_leaves=2,
                                random_state=42, subsample=0.9))])"
$ws.Range("B2").Value = 0.6476190476190476
$ws.Range("C2").Value = "{'selector': None, 'scaler': RobustScaler(), 'model__subsample': 0.9, 'model__num_leaves': 2, 'model__min_child_samples': 1, 'model__max_depth': 7, 'model__learning_rate': 0.01, 'model__colsample_bytree': 0.7, 'model__class_weight': 'balanced', 'model__boosting_type': 'dart'}"
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = "[1 1 0 0 1 0 0 0 0 1 0 1]"
$ws.Range("F2").Value = "[0 0 1 0 0 0 1 0 1 0 0 0]"
$ws.Range("G2").Value = 77
$ws.Range("H2").Value = 0.8636469221835076
$ws.Range("I2").Value = 0.01793771493432178
$ws.Range("J2").Value = 0.5710801393728222
$ws.Range("K2").Value = 0.05990602028042255

# Row 3
$ws.Range("A3").Value = "Pipeline(steps=[('scaler', None),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f45a97add00>),
                ('model',
                 LGBMClassifier(class_weight='balanced', colsample_bytree=0.7,
                                learning_rate=0.01, max_depth=1, num_leaves=10,
                                random_state=42, subsample=0.5))])"
$ws.Range("B3").Value = 0.638095238095238
$ws.Range("C3").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7f3a6c5d4ca0>, 'scaler': None, 'model__subsample': 0.5, 'model__num_leaves': 10, 'model__min_child_samples': 20, 'model__max_depth': 1, 'model__learning_rate': 0.01, 'model__colsample_bytree': 0.7, 'model__class_weight': 'balanced', 'model__boosting_type': 'gbdt'}"
$ws.Range("D3").Value = 0.5
$ws.Range("E3").Value = "[1 1 0 1 0 0 1 0 1 1 1 0]"
$ws.Range("F3").Value = "[0 1 0 1 1 0 0 1 0 1 0 0]"
$ws.Range("G3").Value = 69
$ws.Range("H3").Value = 0.8768796992481203
$ws.Range("I3").Value = 0.01784000748549297
$ws.Range("J3").Value = 0.5106516290726816
$ws.Range("K3").Value = 0.07478733123258723

# Row 4
$ws.Range("A4").Value = "Pipeline(steps=[('scaler', None),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f45a9c158b0>),
                ('model',
                 LGBMClassifier(class_weight='balanced', colsample_bytree=0.5,
                                learning_rate=0.05, max_depth=5, num_leaves=10,
                                random_state=42, subsample=0.7))])"
$ws.Range("B4").Value = 0.6095238095238095
$ws.Range("C4").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7f3a642e0100>, 'scaler': None, 'model__subsample': 0.7, 'model__num_leaves': 10, 'model__min_child_samples': 20, 'model__max_depth': 5, 'model__learning_rate': 0.05, 'model__colsample_bytree': 0.5, 'model__class_weight': 'balanced', 'model__boosting_type': 'gbdt'}"
$ws.Range("D4").Value = 0.888888888888889
$ws.Range("E4").Value = "[1 0 1 1 1 1 0 1 0 1 0 1]"
$ws.Range("F4").Value = "[1 1 1 1 1 1 0 1 1 1 0 1]"
$ws.Range("G4").Value = 42
$ws.Range("H4").Value = 0.874624060150376
$ws.Range("I4").Value = 0.01894050260234806
$ws.Range("J4").Value = 0.5081453634085213
$ws.Range("K4").Value = 0.09603989856371536

# Row 5
$ws.Range("A5").Value = "Pipeline(steps=[('scaler', None),
                ('selector',
                 <__main__.NamedFeatureSelector object at 0x7f3a642e02b0>),
                ('model',
                 LGBMClassifier(boosting_type='dart', class_weight='balanced',
                                colsample_bytree=0.7, learning_rate=0.2,
                                max_depth=1, num_leaves=2, random_state=42,
                                subsample=0.9))])"
$ws.Range("B5").Value = 0.6476190476190475
$ws.Range("C5").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7f3a643fe9a0>, 'scaler': None, 'model__subsample': 0.9, 'model__num_leaves': 2, 'model__min_child_samples': 20, 'model__max_depth': 1, 'model__learning_rate': 0.2, 'model__colsample_bytree': 0.7, 'model__class_weight': 'balanced', 'model__boosting_type': 'dart'}"
$ws.Range("D5").Value = 0.5714285714285714
$ws.Range("E5").Value = "[1 1 0 0 0 0 1 0 1 1 1 1]"
$ws.Range("F5").Value = "[0 0 0 1 0 1 1 1 1 0 1 1]"
$ws.Range("G5").Value = 11
$ws.Range("H5").Value = 0.8723684210526316
$ws.Range("I5").Value = 0.01862244105769229
$ws.Range("J5").Value = 0.5213032581453634
$ws.Range("K5").Value = 0.06755114146166681

# Row 6
$ws.Range("A6").Value = "Pipeline(steps=[('scaler', MinMaxScaler()),
                ('selector',
                 SelectFromModel(estimator=ExtraTreesClassifier(random_state=42))),
                ('model',
                 LGBMClassifier(colsample_bytree=0.9, learning_rate=0.01,
                                max_depth=1, min_child_samples=10, num_leaves=5,
                                random_state=42, subsample=0.9))])"
$ws.Range("B6").Value = 0.6
$ws.Range("C6").Value = "{'selector': SelectFromModel(estimator=ExtraTreesClassifier(random_state=42)), 'scaler': MinMaxScaler(), 'model__subsample': 0.9, 'model__num_leaves': 5, 'model__min_child_samples': 10, 'model__max_depth': 1, 'model__learning_rate': 0.01, 'model__colsample_bytree': 0.9, 'model__class_weight': None, 'model__boosting_type': 'gbdt'}"
$ws.Range("D6").Value = 0.5
$ws.Range("E6").Value = "[1 1 1 1 0 0 0 0 1 1 0 0]"
$ws.Range("F6").Value = "[1 1 1 0 1 1 1 1 1 0 1 1]"
$ws.Range("G6").Value = 14
$ws.Range("H6").Value = 0.8876825396825396
$ws.Range("I6").Value = 0.0167951775428383
$ws.Range("J6").Value = 0.5179682539682539
$ws.Range("K6").Value = 0.08350895700281397
